$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match the repulled/recalculated data
$ws.Range("F2").Value = 7
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -8
$ws.Range("F12").Value = 1
